# Add a new "FAVÖK" (EBITDA) row to the DATA_RAW sheet (row 23) and to the
# "gelir tablosu (çeyreklik)" sheet (row 9), with the same set of values.

$wb = $excel.ActiveWorkbook

$values = @(0,0,0,0,0,0,217000000,154000000,0,244000000,159000000,251000000,0,394000000,270000000,268000000,0,322000000,201000000,79000000,405000000,192000000,107000000,-28000000,510000000,181000000,79000000,-25000000,591000000,146000000,-15000000,-544000000,1108000000,255000000,464000000,406000000,1255000000,1071000000,3194000000,2659000000,2127000000,3440000000,3261000000,3259000000,3786000000,5598000000,6613000000,7502000000,3879000000,9274000000,8286000000,8720000000,6873000000,9351000000,11591000000,12717000000,0,16373000000,15740000000,13800000000,0,12862000000,13033000000,12852000000,0,15897000000,19270000000,23319000000,0,26991000000,26710000000,30853000000,32667000000,34398000000,34218000000)

# --- DATA_RAW sheet: new row 23 ---
$wsData = $wb.Worksheets.Item("DATA_RAW")
$wsData.Cells.Item(23, 1).Value = "FAVÖK"
for ($i = 0; $i -lt $values.Length; $i++) {
    $wsData.Cells.Item(23, $i + 2).Value = $values[$i]
}

# --- "gelir tablosu (çeyreklik)" sheet: new row 9 ---
$wsGelir = $wb.Worksheets.Item("gelir tablosu (çeyreklik)")
$wsGelir.Cells.Item(9, 1).Value = "FAVÖK"
for ($i = 0; $i -lt $values.Length; $i++) {
    $wsGelir.Cells.Item(9, $i + 2).Value = $values[$i]
}
